$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = 50
$ws.Range("B11").Select() | Out-Null
